$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeXmlName {
    param(
        $Shape,
        [string]$OldName,
        [string]$NewName
    )
    $rng = $Shape.Range
    $xml = $rng.WordOpenXML
    $search = 'name="' + $OldName + '"'
    $replace = 'name="' + $NewName + '"'
    $newXml = $xml.Replace($search, $replace)
    $rng.InsertXML($newXml)
}

# First-page header (header1.xml): BTec logo "image1.jpg" -> "image2.jpg"
$hdrFirst = $sec.Headers.Item(2)
$hdrShape = $hdrFirst.Range.InlineShapes.Item(1)
Rename-InlineShapeXmlName -Shape $hdrShape -OldName "image1.jpg" -NewName "image2.jpg"

# Primary footer (footer2.xml): Pearson logo "image2.png" -> "image1.png"
$ftrPrimary = $sec.Footers.Item(1)
$ftrPrimaryShape = $ftrPrimary.Range.InlineShapes.Item(1)
Rename-InlineShapeXmlName -Shape $ftrPrimaryShape -OldName "image2.png" -NewName "image1.png"

# First-page footer (footer1.xml): Pearson logo "image2.png" -> "image1.png"
$ftrFirst = $sec.Footers.Item(2)
$ftrFirstShape = $ftrFirst.Range.InlineShapes.Item(1)
Rename-InlineShapeXmlName -Shape $ftrFirstShape -OldName "image2.png" -NewName "image1.png"

Write-Output "done"
